$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.627.10"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.292.01"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.50"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.77"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.52"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.78"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "2.636.43"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.847"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "2.299.09"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "43.616.86"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.80"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  +5.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.47"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  -8.85%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.47"
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.03"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.41"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.87"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.236"
$ws.Range("E40").Value = "  -5.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.29"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  +11.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.19"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.23"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "2.514.98"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.428"
$ws.Range("E51").Value = "  -0.66%  "
